# Update the cryptos price snapshot (columns D "Price" and E "Volume(1h)").
# Values are stored as plain text strings. For the subset of new Price values
# that are valid-looking numbers (e.g. "1.00", "0.600"), force the cell to Text
# format first so Excel keeps the exact original text instead of normalizing it
# into a number (dropping trailing zeros, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.907.78"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "2.670.03"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.49"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.58"
$ws.Range("E9").Value = "  +1.33%  "

$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +4.42%  "

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").Value = "3.143.97"
$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.73"
$ws.Range("E14").Value = "  +10.77%  "

$ws.Range("D15").Value = "60.907.15"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("D17").Value = "2.671.31"
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.58"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("E19").Value = "  +1.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.79"
$ws.Range("E20").Value = "  +0.93%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("E23").Value = "  +1.40%  "

$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("E26").Value = "  +1.43%  "

$ws.Range("E27").Value = "  +4.82%  "

$ws.Range("E28").Value = "  +7.24%  "

$ws.Range("E29").Value = "  +2.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.82"
$ws.Range("E30").Value = "  +6.86%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.43"
$ws.Range("E32").Value = "  +2.88%  "

$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("E34").Value = "  +8.76%  "

$ws.Range("E35").Value = "  +5.41%  "

$ws.Range("E36").Value = "  +7.66%  "

$ws.Range("E37").Value = "  +3.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "332.05"
$ws.Range("E38").Value = "  +12.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.51"
$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.880"
$ws.Range("E41").Value = "  +4.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.58"
$ws.Range("E42").Value = "  +3.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.80"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("E45").Value = "  +1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0562"
$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.616"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("E48").Value = "  +3.33%  "

$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("E50").Value = "  +3.07%  "

$ws.Range("D51").Value = "2.115.99"
$ws.Range("E51").Value = "  +4.48%  "
